$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# New cells (no prior value) need the "h:mm" time number format applied
# so they reuse the same style as the other B-column timing cells.
$newCells = @("B9", "B12", "B14", "B15", "B18", "B19", "B23")
foreach ($addr in $newCells) {
    $ws.Range($addr).NumberFormat = "h:mm"
}

$ws.Range("B8").Value  = 0.52500000000000002
$ws.Range("B9").Value  = 0.52638888888888891
$ws.Range("B12").Value = 0.53611111111111109
$ws.Range("B14").Value = 0.54375000000000007
$ws.Range("B15").Value = 0.5444444444444444
$ws.Range("B18").Value = 0.5541666666666667
$ws.Range("B19").Value = 0.5541666666666667
$ws.Range("B23").Value = 0.55972222222222223

# Move the active selection (also clears the stale topLeftCell scroll state)
$ws.Range("B24").Select()
